$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Use case 2 (row 8, A8=7) - B8 should now contain the number 1 instead of the text "user7"
$ws.Range("B8").Value = 1

# Update the selection to reflect the new active cell B8
$ws.Range("B8").Select()
